$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 305, shifting existing rows 305:366 down to 306:367
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row 305 with its data
$ws.Range("A305").Value = 6
$ws.Range("B305").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C305").Value = "Metropolitana"
$ws.Range("D305").Value = 45211
$ws.Range("D305").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E305").Value = 13
$ws.Range("F305").Value = 100112029
$ws.Range("G305").Value = "Orégano"
$ws.Range("H305").Value = "Sin especificar"
$ws.Range("I305").Value = "Primera"
$ws.Range("J305").Value = 29
$ws.Range("K305").Value = 16000
$ws.Range("L305").Value = 16000
$ws.Range("M305").Value = 16000
$ws.Range("N305").Value = "`$/docena de atados"
$ws.Range("O305").Value = "Región Metropolitana"
$ws.Range("P305").Value = 5333
$ws.Range("Q305").Value = 3
$ws.Range("R305").Value = "Hortaliza"
